$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values, replacing the previous "Strike#" values.
$kValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 4
    20 = 2
    21 = 2
    22 = 1
    23 = 2
    24 = 2
    25 = 1
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
